$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 9 (C9:I9) keeps the same formula "=100*$C$5*$C$4" it already had; the
# diff only moves which cell anchors the shared-formula group, which has no
# effect on the evaluated values, so nothing needs to change here.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# New block starting at row 16: cash-flow numbers 1..8
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 6
$ws.Range("I16").Value = 7
$ws.Range("J16").Value = 8

# Row 17: coupon amount formulas, mirroring row 9, but J17 stays a literal "=4"
$ws.Range("C17").Formula = "=100*`$C`$5*`$C`$4"
$ws.Range("D17:I17").Formula = "=100*`$C`$5*`$C`$4"
$ws.Range("J17").Formula = "=4"

# Row 18: discount-factor formulas, mirroring row 10, with the 0.00000 number format
$ws.Range("C18:J18").Formula = "=EXP(-`$C`$6*C16*`$C`$4)"

# Row 20: present value of each cash flow, mirroring row 12
# (computed before the row-18 number format is applied, so the format used on
# row 18 does not leak onto the dependent formulas in row 20)
$ws.Range("C20").Formula = "=C17*C18"
$ws.Range("D20:J20").Formula = "=D17*D18"

# Row 22: running totals (I22 intentionally left blank, matching source)
$ws.Range("F22").Formula = "=SUM(C20:F20)"
$ws.Range("G22").Formula = "=SUM(C20:G20)"
$ws.Range("H22").Formula = "=SUM(C20:H20)"
$ws.Range("J22").Formula = "=SUM(C20:J20)"

# Now apply the discount-factor number format to row 18
$ws.Range("C18:J18").NumberFormat = "0.00000"

# Row 25: simple scratch calculation
$ws.Range("C25").Formula = "=100*0.08*0.5"

# Column J width adjustment
$ws.Columns(10).ColumnWidth = 11

# Move the active selection to the new work area, like in the edited file
$ws.Range("E16").Select()
